$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14, shifting the existing rows 14-17 down to 15-18.
$ws.Rows("14:14").Insert()

# Populate the new row 14 with the weekly price-report entry.
$ws.Range("A14").Value = 8
$ws.Range("B14").Value = "Terminal La Palmera de La Serena"
$ws.Range("C14").Value = "Coquimbo"
$ws.Range("D14").Value = 44855
$ws.Range("E14").Value = 4
$ws.Range("F14").Value = 100112013
$ws.Range("G14").Value = "Alcachofa"
$ws.Range("H14").Value = "Española"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 540
$ws.Range("K14").Value = 9500
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = 9750
$ws.Range("N14").Value = "$/caja 30 unidades"
$ws.Range("O14").Value = "Provincia del Elquí"
$ws.Range("P14").Value = 325
$ws.Range("Q14").Value = 30
$ws.Range("R14").Value = "Hortaliza"
